$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 839.625
$ws.Cells.Item(98, 9).Value = 839.625
$ws.Cells.Item(98, 11).Value = 839.625
$ws.Cells.Item(98, 13).Value = 658.375
$ws.Cells.Item(122, 8).Value = 839.625
$ws.Cells.Item(122, 9).Value = 839.625
$ws.Cells.Item(122, 11).Value = 2518.875
$ws.Cells.Item(122, 13).Value = -68.875
$ws.Cells.Item(137, 8).Value = 1480.4286
$ws.Cells.Item(137, 9).Value = 837.7
$ws.Cells.Item(137, 10).Value = 3087.25
$ws.Cells.Item(137, 11).Value = 2513.1
$ws.Cells.Item(137, 12).Value = 9261.75
$ws.Cells.Item(137, 13).Value = 36.89999999999964
$ws.Cells.Item(137, 14).Value = -14361.75
$ws.Cells.Item(138, 8).Value = 4142.923
$ws.Cells.Item(138, 9).Value = 1182.25
$ws.Cells.Item(138, 10).Value = 4813.264
$ws.Cells.Item(138, 11).Value = 3546.75
$ws.Cells.Item(138, 12).Value = 14439.792
$ws.Cells.Item(138, 13).Value = 1593.25
$ws.Cells.Item(138, 14).Value = -24719.792

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3225.913
$ws.Cells.Item(61, 9).Value = 1593.7894
$ws.Cells.Item(61, 11).Value = 1593.7894
$ws.Cells.Item(61, 13).Value = -1381.7894
$ws.Cells.Item(74, 8).Value = 9733.166999999999
$ws.Cells.Item(74, 9).Value = 11378.4
$ws.Cells.Item(74, 11).Value = 11378.4
$ws.Cells.Item(74, 13).Value = -10504.4
$ws.Cells.Item(77, 8).Value = 9733.166999999999
$ws.Cells.Item(77, 9).Value = 11378.4
$ws.Cells.Item(77, 11).Value = 56892
$ws.Cells.Item(77, 13).Value = -52524
$ws.Cells.Item(132, 8).Value = 2611.6667
$ws.Cells.Item(132, 9).Value = 2156.9285
$ws.Cells.Item(132, 10).Value = 4203.25
$ws.Cells.Item(132, 11).Value = 6470.7855
$ws.Cells.Item(132, 12).Value = 12609.75
$ws.Cells.Item(132, 13).Value = -3940.7855
$ws.Cells.Item(132, 14).Value = -17669.75
$ws.Cells.Item(136, 8).Value = 3225.913
$ws.Cells.Item(136, 9).Value = 1593.7894
$ws.Cells.Item(136, 11).Value = 4781.3682
$ws.Cells.Item(136, 13).Value = -2231.3682
$ws.Cells.Item(138, 8).Value = 18709.5
$ws.Cells.Item(138, 10).Value = 18709.5
$ws.Cells.Item(138, 12).Value = 18709.5
$ws.Cells.Item(138, 14).Value = -28989.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(38, 8).Value = 45000
$ws.Cells.Item(38, 10).Value = 45000
$ws.Cells.Item(38, 12).Value = 45000
$ws.Cells.Item(38, 14).Value = -45832
$ws.Cells.Item(86, 8).Value = 1878.8889
$ws.Cells.Item(86, 9).Value = 1840.4
$ws.Cells.Item(86, 10).Value = 2071.3333
$ws.Cells.Item(86, 11).Value = 1840.4
$ws.Cells.Item(86, 12).Value = 2071.3333
$ws.Cells.Item(86, 13).Value = -717.4000000000001
$ws.Cells.Item(86, 14).Value = -4317.3333
$ws.Cells.Item(89, 8).Value = 1878.8889
$ws.Cells.Item(89, 9).Value = 1840.4
$ws.Cells.Item(89, 10).Value = 2071.3333
$ws.Cells.Item(89, 11).Value = 9202
$ws.Cells.Item(89, 12).Value = 10356.6665
$ws.Cells.Item(89, 13).Value = -3586
$ws.Cells.Item(89, 14).Value = -21588.6665
$ws.Cells.Item(99, 8).Value = 1750.3077
$ws.Cells.Item(99, 9).Value = 1317.6
$ws.Cells.Item(99, 11).Value = 1317.6
$ws.Cells.Item(99, 13).Value = 180.4000000000001
$ws.Cells.Item(134, 8).Value = 2056.88
$ws.Cells.Item(134, 9).Value = 1700.3636
$ws.Cells.Item(134, 10).Value = 4671.3335
$ws.Cells.Item(134, 11).Value = 5101.0908
$ws.Cells.Item(134, 12).Value = 14014.0005
$ws.Cells.Item(134, 13).Value = -2566.0908
$ws.Cells.Item(134, 14).Value = -19084.0005
$ws.Cells.Item(140, 8).Value = 55000
$ws.Cells.Item(140, 10).Value = 55000
$ws.Cells.Item(140, 12).Value = -65360

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2078.3264
$ws.Cells.Item(31, 9).Value = 1488.625
$ws.Cells.Item(31, 10).Value = 3188.353
$ws.Cells.Item(31, 11).Value = 1488.625
$ws.Cells.Item(31, 12).Value = 3188.353
$ws.Cells.Item(31, 13).Value = -1193.625
$ws.Cells.Item(31, 14).Value = -3778.353
$ws.Cells.Item(34, 8).Value = 2078.3264
$ws.Cells.Item(34, 9).Value = 1488.625
$ws.Cells.Item(34, 10).Value = 3188.353
$ws.Cells.Item(34, 11).Value = 1488.625
$ws.Cells.Item(34, 12).Value = 3188.353
$ws.Cells.Item(34, 13).Value = -1286.625
$ws.Cells.Item(34, 14).Value = -3592.353
$ws.Cells.Item(58, 8).Value = 1452.0454
$ws.Cells.Item(58, 9).Value = 1283.0952
$ws.Cells.Item(58, 11).Value = 1283.0952
$ws.Cells.Item(58, 13).Value = -1080.0952
$ws.Cells.Item(132, 8).Value = 2636.6206
$ws.Cells.Item(132, 9).Value = 1906
$ws.Cells.Item(132, 10).Value = 3832.182
$ws.Cells.Item(132, 11).Value = 5718
$ws.Cells.Item(132, 12).Value = 11496.546
$ws.Cells.Item(132, 13).Value = -3188
$ws.Cells.Item(132, 14).Value = -16556.546
$ws.Cells.Item(134, 8).Value = 3317.25
$ws.Cells.Item(134, 9).Value = 1821
$ws.Cells.Item(134, 11).Value = 5463
$ws.Cells.Item(134, 13).Value = -2928
$ws.Cells.Item(136, 8).Value = 1452.0454
$ws.Cells.Item(136, 9).Value = 1283.0952
$ws.Cells.Item(136, 11).Value = 3849.2856
$ws.Cells.Item(136, 13).Value = -1299.2856

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 875.03845
$ws.Cells.Item(131, 9).Value = 573.75
$ws.Cells.Item(131, 10).Value = 929.8182
$ws.Cells.Item(131, 11).Value = 1721.25
$ws.Cells.Item(131, 12).Value = 2789.4546
$ws.Cells.Item(131, 13).Value = 3318.75
$ws.Cells.Item(131, 14).Value = -12869.4546
$ws.Cells.Item(137, 8).Value = 3297.2307
$ws.Cells.Item(137, 9).Value = 2685.5715
$ws.Cells.Item(137, 10).Value = 3522.5789
$ws.Cells.Item(137, 11).Value = 8056.7145
$ws.Cells.Item(137, 12).Value = 10567.7367
$ws.Cells.Item(137, 13).Value = -2956.7145
$ws.Cells.Item(137, 14).Value = -20767.7367

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2086.2
$ws.Cells.Item(102, 9).Value = 2038.1052
$ws.Cells.Item(102, 10).Value = 3000
$ws.Cells.Item(102, 11).Value = 2038.1052
$ws.Cells.Item(102, 12).Value = 3000
$ws.Cells.Item(102, 13).Value = -6244
$ws.Cells.Item(122, 8).Value = 25001388
$ws.Cells.Item(122, 9).Value = 25001388
$ws.Cells.Item(122, 11).Value = 75004164
$ws.Cells.Item(122, 13).Value = -75001714
$ws.Cells.Item(132, 8).Value = 5475.08
$ws.Cells.Item(132, 9).Value = 5449.9546
$ws.Cells.Item(132, 10).Value = 5659.3335
$ws.Cells.Item(132, 11).Value = 16349.8638
$ws.Cells.Item(132, 12).Value = 16978.0005
$ws.Cells.Item(132, 13).Value = -13819.8638
$ws.Cells.Item(132, 14).Value = -22038.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 5654.3516
$ws.Cells.Item(132, 9).Value = 6057.0806
$ws.Cells.Item(132, 10).Value = 3573.5833
$ws.Cells.Item(132, 11).Value = 18171.2418
$ws.Cells.Item(132, 12).Value = 10720.7499
$ws.Cells.Item(132, 13).Value = -15641.2418
$ws.Cells.Item(132, 14).Value = -15780.7499
$ws.Cells.Item(135, 8).Value = 14490
$ws.Cells.Item(135, 10).Value = 14490
$ws.Cells.Item(135, 12).Value = 14490
$ws.Cells.Item(135, 14).Value = -24630
$ws.Cells.Item(136, 8).Value = 1749.0571
$ws.Cells.Item(136, 9).Value = 1307.0667
$ws.Cells.Item(136, 10).Value = 4401
$ws.Cells.Item(136, 11).Value = 3921.2001
$ws.Cells.Item(136, 12).Value = 13203
$ws.Cells.Item(136, 13).Value = -1371.2001
$ws.Cells.Item(136, 14).Value = -18303

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 333335840
$ws.Cells.Item(14, 10).Value = 333335840
$ws.Cells.Item(14, 12).Value = -333336176
$ws.Cells.Item(46, 8).Value = 45286
$ws.Cells.Item(46, 10).Value = 45286
$ws.Cells.Item(46, 12).Value = -45748
$ws.Cells.Item(122, 8).Value = 2472.35
$ws.Cells.Item(122, 9).Value = 2117.4614
$ws.Cells.Item(122, 10).Value = 3131.4285
$ws.Cells.Item(122, 11).Value = 6352.3842
$ws.Cells.Item(122, 12).Value = 9394.2855
$ws.Cells.Item(122, 13).Value = -3902.3842
$ws.Cells.Item(122, 14).Value = -14294.2855
$ws.Cells.Item(132, 8).Value = 3619.75
$ws.Cells.Item(132, 9).Value = 3188.5
$ws.Cells.Item(132, 10).Value = 4051
$ws.Cells.Item(132, 11).Value = 9565.5
$ws.Cells.Item(132, 12).Value = 12153
$ws.Cells.Item(132, 13).Value = -7035.5
$ws.Cells.Item(132, 14).Value = -17213
$ws.Cells.Item(134, 8).Value = 45286
$ws.Cells.Item(134, 10).Value = 45286
$ws.Cells.Item(134, 12).Value = -140928
$ws.Cells.Item(136, 8).Value = 4565.3438
$ws.Cells.Item(136, 9).Value = 4485.2144
$ws.Cells.Item(136, 10).Value = 5126.25
$ws.Cells.Item(136, 11).Value = 13455.6432
$ws.Cells.Item(136, 12).Value = 15378.75
$ws.Cells.Item(136, 13).Value = -10905.6432
$ws.Cells.Item(136, 14).Value = -20478.75
